# The sheet contains a weekly price log where each row holds one day's
# record. A new, more recent record is inserted as the new row 3, which
# pushes all the former rows 3-19 down to rows 4-20 (the newest entries
# are kept near the top, right under the header/row2).
#
# Columns: A Mercado ID, B Mercado, C Región, D Fecha, E Codreg,
#          F Categoría ID, G Categoría, H Variedad, I Calidad, J Volumen,
#          K Precio mínimo, L Precio máximo, M Precio promedio ponderado,
#          N Unidad de comercialización, O Origen, P Precio $/Kg,
#          Q Kg o Unidades, R Clasificación

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 3; Excel automatically shifts the
# existing rows 3..19 down to 4..20 and extends the used range/dimension.
$ws.Rows("3:3").Insert()

# Fill the newly inserted row 3. Columns that are constant across every
# record in this sheet (A, B, C, E, F, G, H, I, N, O, Q, R) are copied
# from the row directly below (row 4, which now holds what used to be
# row 3), while the data columns that actually change (D, J, K, L, M, P)
# get the new reported values.
$ws.Cells(3,1).Value2  = $ws.Cells(4,1).Value2    # Mercado ID
$ws.Cells(3,2).Value2  = $ws.Cells(4,2).Value2    # Mercado
$ws.Cells(3,3).Value2  = $ws.Cells(4,3).Value2    # Región
$ws.Cells(3,4).Value2  = 44882                    # Fecha
$ws.Cells(3,5).Value2  = $ws.Cells(4,5).Value2    # Codreg
$ws.Cells(3,6).Value2  = $ws.Cells(4,6).Value2    # Categoría ID
$ws.Cells(3,7).Value2  = $ws.Cells(4,7).Value2    # Categoría
$ws.Cells(3,8).Value2  = $ws.Cells(4,8).Value2    # Variedad
$ws.Cells(3,9).Value2  = $ws.Cells(4,9).Value2    # Calidad
$ws.Cells(3,10).Value2 = 70                       # Volumen
$ws.Cells(3,11).Value2 = 7000                     # Precio mínimo
$ws.Cells(3,12).Value2 = 7000                     # Precio máximo
$ws.Cells(3,13).Value2 = 7000                     # Precio promedio ponderado
$ws.Cells(3,14).Value2 = $ws.Cells(4,14).Value2   # Unidad de comercialización
$ws.Cells(3,15).Value2 = $ws.Cells(4,15).Value2   # Origen
$ws.Cells(3,16).Value2 = 438                      # Precio $/Kg
$ws.Cells(3,17).Value2 = $ws.Cells(4,17).Value2   # Kg o Unidades
$ws.Cells(3,18).Value2 = $ws.Cells(4,18).Value2   # Clasificación
